# Auto-generated cell value updates applying the Moogle_Profits market-price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2929.2727
$ws.Range("I2").Value = 1644.5
$ws.Range("J2").Value = 3214.7778
$ws.Range("K2").Value = 1644.5
$ws.Range("L2").Value = 3214.7778
$ws.Range("M2").Value = -1531.5
$ws.Range("N2").Value = -3440.7778
$ws.Range("H17").Value = 5308.625
$ws.Range("J17").Value = 5308.625
$ws.Range("L17").Value = 15925.875
$ws.Range("N17").Value = -16261.875
$ws.Range("H55").Value = 336
$ws.Range("I55").Value = 207.5
$ws.Range("J55").Value = 464.5
$ws.Range("K55").Value = 207.5
$ws.Range("L55").Value = 464.5
$ws.Range("M55").Value = 6.5
$ws.Range("N55").Value = -892.5
$ws.Range("H70").Value = 7219.222
$ws.Range("J70").Value = 9011.714
$ws.Range("L70").Value = 27035.142
$ws.Range("N70").Value = -27575.142
$ws.Range("H73").Value = 7219.222
$ws.Range("J73").Value = 9011.714
$ws.Range("L73").Value = 27035.142
$ws.Range("N73").Value = -28907.142
$ws.Range("H101").Value = 699.2
$ws.Range("I101").Value = 249.5
$ws.Range("J101").Value = 2498
$ws.Range("K101").Value = 748.5
$ws.Range("L101").Value = 7494
$ws.Range("M101").Value = 873.5
$ws.Range("N101").Value = -10738
$ws.Range("H113").Value = 2194.9534
$ws.Range("J113").Value = 2919.7778
$ws.Range("L113").Value = 2919.7778
$ws.Range("N113").Value = -9427.7778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 642.2
$ws.Range("I2").Value = 642.2
$ws.Range("K2").Value = 642.2
$ws.Range("M2").Value = -529.2
$ws.Range("H32").Value = 5879.9844
$ws.Range("I32").Value = 2764.3965
$ws.Range("K32").Value = 2764.3965
$ws.Range("M32").Value = -2477.3965
$ws.Range("H45").Value = 2715.8696
$ws.Range("I45").Value = 2366.7368
$ws.Range("K45").Value = 2366.7368
$ws.Range("M45").Value = -1989.7368
$ws.Range("H102").Value = 1505
$ws.Range("I102").Value = 1505
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1505
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 117
$ws.Range("N102").ClearContents()
$ws.Range("H116").Value = 642.2
$ws.Range("I116").Value = 642.2
$ws.Range("K116").Value = 642.2
$ws.Range("M116").Value = 1651.8
$ws.Range("H124").Value = 19166.666
$ws.Range("J124").Value = 19166.666
$ws.Range("L124").Value = 19166.666
$ws.Range("N124").Value = -28986.666
$ws.Range("H132").Value = 1706.6444
$ws.Range("I132").Value = 1163.8379
$ws.Range("J132").Value = 4217.125
$ws.Range("K132").Value = 3491.5137
$ws.Range("L132").Value = 12651.375
$ws.Range("M132").Value = -961.5137
$ws.Range("N132").Value = -17711.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 642.2
$ws.Range("I3").Value = 642.2
$ws.Range("K3").Value = 642.2
$ws.Range("M3").Value = -528.2
$ws.Range("H20").Value = 9239.883
$ws.Range("I20").Value = 2250.6667
$ws.Range("J20").Value = 17102.75
$ws.Range("K20").Value = 2250.6667
$ws.Range("L20").Value = 17102.75
$ws.Range("M20").Value = -2003.6667
$ws.Range("N20").Value = -17596.75
$ws.Range("H80").Value = 40055.2
$ws.Range("J80").Value = 69
$ws.Range("L80").Value = 69
$ws.Range("N80").Value = -2065
$ws.Range("H83").Value = 40055.2
$ws.Range("J83").Value = 69
$ws.Range("L83").Value = 345
$ws.Range("N83").Value = -10329
$ws.Range("H105").Value = 4343.222
$ws.Range("I105").Value = 3916.6667
$ws.Range("K105").Value = 3916.6667
$ws.Range("M105").Value = -2169.6667
$ws.Range("H134").Value = 4263.676
$ws.Range("I134").Value = 3180.5454
$ws.Range("K134").Value = 9541.636200000001
$ws.Range("M134").Value = -7006.636200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 885
$ws.Range("I16").Value = 377.5
$ws.Range("K16").Value = 377.5
$ws.Range("M16").Value = -90.5
$ws.Range("H62").Value = 6941.7646
$ws.Range("I62").Value = 5236.3335
$ws.Range("J62").Value = 7872
$ws.Range("K62").Value = 5236.3335
$ws.Range("L62").Value = 7872
$ws.Range("M62").Value = -4612.3335
$ws.Range("N62").Value = -9120
$ws.Range("H65").Value = 6941.7646
$ws.Range("I65").Value = 5236.3335
$ws.Range("J65").Value = 7872
$ws.Range("K65").Value = 26181.6675
$ws.Range("L65").Value = 39360
$ws.Range("M65").Value = -23061.6675
$ws.Range("N65").Value = -45600
$ws.Range("H99").Value = 4347.1816
$ws.Range("I99").Value = 4019.4285
$ws.Range("J99").Value = 4920.75
$ws.Range("K99").Value = 4019.4285
$ws.Range("L99").Value = 4920.75
$ws.Range("M99").Value = -2521.4285
$ws.Range("N99").Value = -7916.75
$ws.Range("H113").Value = 885
$ws.Range("I113").Value = 377.5
$ws.Range("K113").Value = 377.5
$ws.Range("M113").Value = 1792.5
$ws.Range("H122").Value = 3187.3333
$ws.Range("I122").Value = 2316.4167
$ws.Range("K122").Value = 6949.250100000001
$ws.Range("M122").Value = -4499.250100000001
$ws.Range("H126").Value = 4347.1816
$ws.Range("I126").Value = 4019.4285
$ws.Range("J126").Value = 4920.75
$ws.Range("K126").Value = 12058.2855
$ws.Range("L126").Value = 14762.25
$ws.Range("M126").Value = -9588.2855
$ws.Range("N126").Value = -19702.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 185.33333
$ws.Range("I2").Value = 9.6
$ws.Range("J2").Value = 405
$ws.Range("K2").Value = 57.59999999999999
$ws.Range("L2").Value = 2430
$ws.Range("M2").Value = 55.40000000000001
$ws.Range("N2").Value = -2656
$ws.Range("H140").Value = 2170.9768
$ws.Range("I140").Value = 1883.875
$ws.Range("K140").Value = 5651.625
$ws.Range("M140").Value = -471.625
$ws.Range("H141").Value = 7897.2
$ws.Range("I141").Value = 6023.75
$ws.Range("K141").Value = 18071.25
$ws.Range("M141").Value = -12891.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6956
$ws.Range("I70").Value = 7946.3335
$ws.Range("J70").Value = 3985
$ws.Range("K70").Value = 7946.3335
$ws.Range("L70").Value = 3985
$ws.Range("M70").Value = -7676.3335
$ws.Range("N70").Value = -4525
$ws.Range("H73").Value = 6956
$ws.Range("I73").Value = 7946.3335
$ws.Range("J73").Value = 3985
$ws.Range("K73").Value = 7946.3335
$ws.Range("L73").Value = 3985
$ws.Range("M73").Value = -7010.3335
$ws.Range("N73").Value = -5857
$ws.Range("H80").Value = 2568
$ws.Range("I80").Value = 1602
$ws.Range("J80").Value = 4500
$ws.Range("K80").Value = 1602
$ws.Range("L80").Value = 4500
$ws.Range("M80").Value = -604
$ws.Range("N80").Value = -6496
$ws.Range("H83").Value = 2568
$ws.Range("I83").Value = 1602
$ws.Range("J83").Value = 4500
$ws.Range("K83").Value = 8010
$ws.Range("L83").Value = 22500
$ws.Range("M83").Value = -3018
$ws.Range("N83").Value = -32484
$ws.Range("H122").Value = 4643.3105
$ws.Range("I122").Value = 2232.6875
$ws.Range("J122").Value = 7610.231
$ws.Range("K122").Value = 6698.0625
$ws.Range("L122").Value = 22830.693
$ws.Range("M122").Value = -4248.0625
$ws.Range("N122").Value = -27730.693
$ws.Range("H126").Value = 7220.467
$ws.Range("I126").Value = 7310.364
$ws.Range("K126").Value = 21931.092
$ws.Range("M126").Value = -19461.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2228.2334
$ws.Range("J46").Value = 3122.111
$ws.Range("L46").Value = 3122.111
$ws.Range("N46").Value = -3498.111
$ws.Range("H68").Value = 5997.64
$ws.Range("I68").Value = 3885.8125
$ws.Range("J68").Value = 9752
$ws.Range("K68").Value = 3885.8125
$ws.Range("L68").Value = 9752
$ws.Range("M68").Value = -3136.8125
$ws.Range("N68").Value = -11250
$ws.Range("H71").Value = 5997.64
$ws.Range("I71").Value = 3885.8125
$ws.Range("J71").Value = 9752
$ws.Range("K71").Value = 19429.0625
$ws.Range("L71").Value = 48760
$ws.Range("M71").Value = -15685.0625
$ws.Range("N71").Value = -56248
$ws.Range("H82").Value = 1622.3684
$ws.Range("I82").Value = 813.3333
$ws.Range("J82").Value = 3009.2856
$ws.Range("K82").Value = 813.3333
$ws.Range("L82").Value = 3009.2856
$ws.Range("M82").Value = -452.3333
$ws.Range("N82").Value = -3731.2856
$ws.Range("H85").Value = 1622.3684
$ws.Range("I85").Value = 813.3333
$ws.Range("J85").Value = 3009.2856
$ws.Range("K85").Value = 813.3333
$ws.Range("L85").Value = 3009.2856
$ws.Range("M85").Value = 434.6667
$ws.Range("N85").Value = -5505.2856
$ws.Range("H115").Value = 79999
$ws.Range("J115").Value = 79999
$ws.Range("L115").Value = 79999
$ws.Range("N115").Value = -82349

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 25257.143
$ws.Range("J4").Value = 23560
$ws.Range("L4").Value = 23560
$ws.Range("N4").Value = -23786
$ws.Range("H100").Value = 1090.8235
$ws.Range("I100").Value = 811.125
$ws.Range("K100").Value = 1622.25
$ws.Range("M100").Value = -1081.25
$ws.Range("H110").Value = 94995
$ws.Range("J110").Value = 94995
$ws.Range("L110").Value = 94995
$ws.Range("N110").Value = -103175
